# Revert "[15] Use the same column names as the Jira export."
# Restore the Jira-style (no-space) header names in row 1 and re-add the
# Forecast / Done helper columns (I:L) with their formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprints")

# --- Row 1 headers: revert to Jira-export-style names (no spaces) ---
$ws.Range("A1").Value = "SprintName"
$ws.Range("B1").Value = "StartDate"
$ws.Range("C1").Value = "EndDate"
$ws.Range("D1").Value = "CapacityForecast"
$ws.Range("E1").Value = "EffortForecast"
$ws.Range("F1").Value = "CapacityDone"
$ws.Range("G1").Value = "EffortDone"

# --- Re-add the Forecast / Done computed columns ---
$ws.Range("I1").Value = "Forecast"
$ws.Range("K1").Value = "Done"

$ws.Range("I2").Formula = '=IF(D2<>"",E2/D2,"")'
$ws.Range("J2").Formula = '=I2*D2'
$ws.Range("K2").Formula = '=IF(F2<>"",G2/F2,"")'
$ws.Range("L2").Formula = '=K2*F2'

$ws.Range("I3").Formula = '=IF(D3<>"",E3/D3,"")'
$ws.Range("J3").Formula = '=I3*D3'
$ws.Range("K3").Formula = '=IF(F3<>"",G3/F3,"")'
$ws.Range("L3").Formula = '=K3*F3'

# --- Rows 4-9: restore the older sprint schedule / numbers ---
$ws.Range("F5").Value2 = 3
$ws.Range("G5").ClearContents()

$ws.Range("C6").Value2 = 41760
$ws.Range("D6").Value2 = 6
$ws.Range("E6").ClearContents()

$ws.Range("B7").Value2 = 41761
$ws.Range("C7").Value2 = 41761
$ws.Range("D7").Value2 = 4

$ws.Range("B8").Value2 = 41761
$ws.Range("C8").Value2 = 41761
$ws.Range("D8").Value2 = 3

$ws.Range("B9").Value2 = 41761
$ws.Range("C9").Value2 = 41761
$ws.Range("D9").Value2 = 3

# --- Selection moves from H1:L6 to D10 ---
$ws.Range("D10").Select()
